$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H4").Value = 700
$ws.Range("I4").Value = 700
$ws.Range("K4").Value = 700
$ws.Range("M4").Value = -586

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H43").Value = 1398
$ws.Range("I43").Value = 0
$ws.Range("J43").Value = 1398
$ws.Range("K43").Value = 0
$ws.Range("L43").Value = 1398
$ws.Range("M43").ClearContents()
$ws.Range("N43").Value = -1536

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H70").Value = 15856.143
$ws.Range("I70").Value = 1500
$ws.Range("J70").Value = 18248.834
$ws.Range("K70").Value = 4500
$ws.Range("L70").Value = 54746.50199999999
$ws.Range("M70").Value = -4230
$ws.Range("N70").Value = -55286.50199999999

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H73").Value = 15856.143
$ws.Range("I73").Value = 1500
$ws.Range("J73").Value = 18248.834
$ws.Range("K73").Value = 4500
$ws.Range("L73").Value = 54746.50199999999
$ws.Range("M73").Value = -3564
$ws.Range("N73").Value = -56618.50199999999

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H131").Value = 1598.1666
$ws.Range("J131").Value = 3266.6667
$ws.Range("L131").Value = 9800.000100000001
$ws.Range("N131").Value = -19880.0001

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H132").Value = 926.2157
$ws.Range("I132").Value = 821.7560999999999
$ws.Range("K132").Value = 2465.2683
$ws.Range("M132").Value = 64.73170000000027

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H135").Value = 490.5
$ws.Range("I135").Value = 490.5
$ws.Range("J135").Value = 0
$ws.Range("K135").Value = 4414.5
$ws.Range("L135").Value = 0
$ws.Range("M135").Value = -1879.5
$ws.Range("N135").ClearContents()

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H137").Value = 1973.8462
$ws.Range("I137").Value = 1766
$ws.Range("J137").Value = 2216.3333
$ws.Range("K137").Value = 5298
$ws.Range("L137").Value = 6648.999899999999
$ws.Range("M137").Value = -2748
$ws.Range("N137").Value = -11748.9999

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H97").Value = 1184.2222
$ws.Range("I97").Value = 941.1667
$ws.Range("J97").Value = 1670.3334
$ws.Range("K97").Value = 941.1667
$ws.Range("L97").Value = 1670.3334
$ws.Range("M97").Value = -445.1667
$ws.Range("N97").Value = -2662.3334

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H110").Value = 1282.7931
$ws.Range("I110").Value = 1005.11536
$ws.Range("K110").Value = 1005.11536
$ws.Range("M110").Value = 1039.88464

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H122").Value = 1500
$ws.Range("I122").Value = 1500
$ws.Range("K122").Value = 4500
$ws.Range("M122").Value = -2050

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H132").Value = 1697.6578
$ws.Range("I132").Value = 1172.8966
$ws.Range("K132").Value = 3518.6898
$ws.Range("M132").Value = -988.6898000000001

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 6416.148
$ws.Range("I134").Value = 7419.8184
$ws.Range("K134").Value = 22259.4552
$ws.Range("M134").Value = -19724.4552

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1449.8334
$ws.Range("I16").Value = 1066.6666
$ws.Range("K16").Value = 1066.6666
$ws.Range("M16").Value = -779.6666

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H99").Value = 2600
$ws.Range("I99").Value = 2066.6667
$ws.Range("K99").Value = 2066.6667
$ws.Range("M99").Value = -568.6667000000002

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H105").Value = 2079.889
$ws.Range("I105").Value = 1953.1666
$ws.Range("K105").Value = 1953.1666
$ws.Range("M105").Value = -206.1666

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H113").Value = 1449.8334
$ws.Range("I113").Value = 1066.6666
$ws.Range("K113").Value = 1066.6666
$ws.Range("M113").Value = 1103.3334

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H122").Value = 7338
$ws.Range("I122").Value = 9000
$ws.Range("K122").Value = 27000
$ws.Range("M122").Value = -24550

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H126").Value = 2600
$ws.Range("I126").Value = 2066.6667
$ws.Range("K126").Value = 6200.000100000001
$ws.Range("M126").Value = -3730.000100000001

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H132").Value = 2220.5806
$ws.Range("I132").Value = 1306
$ws.Range("J132").Value = 4141.2
$ws.Range("K132").Value = 3918
$ws.Range("L132").Value = 12423.6
$ws.Range("M132").Value = -1388
$ws.Range("N132").Value = -17483.6

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H3").Value = 2388.7778
$ws.Range("I3").Value = 699.8
$ws.Range("K3").Value = 2099.4
$ws.Range("M3").Value = -1987.4

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H92").Value = 333.66666
$ws.Range("J92").Value = 340.4
$ws.Range("L92").Value = 1021.2
$ws.Range("N92").Value = -3517.2

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H107").Value = 937.4
$ws.Range("J107").Value = 937.4
$ws.Range("L107").Value = 2812.2
$ws.Range("N107").Value = -6652.2

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H113").Value = 17581.834
$ws.Range("J113").Value = 1122.5
$ws.Range("L113").Value = 3367.5
$ws.Range("N113").Value = -7707.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 794.67
$ws.Range("J131").Value = 815.3261
$ws.Range("L131").Value = 2445.9783
$ws.Range("N131").Value = -12525.9783

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H138").Value = 3062.2222
$ws.Range("I138").Value = 2712
$ws.Range("J138").Value = 3500
$ws.Range("K138").Value = 8136
$ws.Range("L138").Value = 10500
$ws.Range("M138").Value = -2996
$ws.Range("N138").Value = -20780

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H113").Value = 929.3570999999999
$ws.Range("I113").Value = 671
$ws.Range("K113").Value = 671
$ws.Range("M113").Value = 1499

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H127").Value = 33088.43
$ws.Range("J127").Value = 33088.43
$ws.Range("L127").Value = 33088.43
$ws.Range("N127").Value = -43008.43

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H139").Value = 58883.332
$ws.Range("J139").Value = 58883.332
$ws.Range("L139").Value = 58883.332
$ws.Range("N139").Value = -69163.33199999999

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 4462.5835
$ws.Range("I7").Value = 2767.4285
$ws.Range("J7").Value = 6835.8
$ws.Range("K7").Value = 2767.4285
$ws.Range("L7").Value = 6835.8
$ws.Range("M7").Value = -2655.4285
$ws.Range("N7").Value = -7059.8

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 4166.3335
$ws.Range("I61").Value = 3749.5
$ws.Range("K61").Value = 3749.5
$ws.Range("M61").Value = -3547.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H93").Value = 642.1667
$ws.Range("I93").Value = 610.6
$ws.Range("K93").Value = 610.6
$ws.Range("M93").Value = 637.4

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H113").Value = 4166.3335
$ws.Range("I113").Value = 3749.5
$ws.Range("K113").Value = 3749.5
$ws.Range("M113").Value = -1579.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H126").Value = 4462.5835
$ws.Range("I126").Value = 2767.4285
$ws.Range("J126").Value = 6835.8
$ws.Range("K126").Value = 8302.2855
$ws.Range("L126").Value = 20507.4
$ws.Range("M126").Value = -5832.2855
$ws.Range("N126").Value = -25447.4

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 2097.2
$ws.Range("I132").Value = 1829.5555
$ws.Range("K132").Value = 5488.666499999999
$ws.Range("M132").Value = -2958.666499999999

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H14").Value = 1086.8
$ws.Range("J14").Value = 1086.8
$ws.Range("L14").Value = 1086.8
$ws.Range("N14").Value = -1422.8

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H69").Value = 13034
$ws.Range("J69").Value = 13034
$ws.Range("L69").Value = 13034
$ws.Range("N69").Value = -14532

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H70").Value = 29500
$ws.Range("J70").Value = 29500
$ws.Range("L70").Value = 29500
$ws.Range("N70").Value = -30130

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H72").Value = 13034
$ws.Range("J72").Value = 13034
$ws.Range("L72").Value = 39102
$ws.Range("N72").Value = -46590

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H73").Value = 29500
$ws.Range("J73").Value = 29500
$ws.Range("L73").Value = 29500
$ws.Range("N73").Value = -31684

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 25255238
$ws.Range("I136").Value = 34724452
$ws.Range("J136").Value = 3999.8333
$ws.Range("K136").Value = 104173356
$ws.Range("L136").Value = 11999.4999
$ws.Range("M136").Value = -104170806
$ws.Range("N136").Value = -17099.4999

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H138").Value = 77776
$ws.Range("J138").Value = 77776
$ws.Range("L138").Value = 77776
$ws.Range("N138").Value = -88056
